$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price records appended to the "Locoto" Hortaliza subconjunto
# (Agrícola del Norte S.A. de Arica) — rows 186 and 187.

$newRows = @(
    @{ Row = 186; A = 1; B = "Agrícola del Norte S.A. de Arica"; C = "Arica y Parinacota"; D = 45191; E = 15; F = 100112042; G = "Locoto"; H = "Sin especificar"; I = "Primera"; J = 230; K = 9000; L = 10000; M = 9652; N = "`$/caja 20 kilos"; O = "Región de Arica y Parinacota"; P = 483; Q = 20; R = "Hortaliza" },
    @{ Row = 187; A = 1; B = "Agrícola del Norte S.A. de Arica"; C = "Arica y Parinacota"; D = 45191; E = 15; F = 100112042; G = "Locoto"; H = "Sin especificar"; I = "Segunda"; J = 190; K = 7000; L = 8000; M = 7526; N = "`$/caja 20 kilos"; O = "Región de Arica y Parinacota"; P = 376; Q = 20; R = "Hortaliza" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
}
